# Applies the textual edits from the commit "Presentation and Risk Assesment"
# to the two Risk Assessment tables in the document.
#
# NOTE 1: Find.Execute in this runtime always searches the whole document body
# regardless of which Range it is invoked on, so every search string below is
# chosen to be unique within the document, and wdReplaceAll (2) is safe to use.
#
# NOTE 2: Table object references returned by $d.Tables.Item(N) are not
# independent - fetching a different table index mutates the previously
# fetched table reference(s) as well, causing stale/incorrect cell lookups.
# To avoid this, each table is (re-)fetched with $d.Tables.Item(N) immediately
# before every cell access, and table 1 and table 2 accesses are never
# interleaved.

$d = $word.ActiveDocument

# --- Table 1, Row 2 (Hacking) ---------------------------------------------

# Risk Statement (col 2): "Using malicious scripts" -> "...scripts trawl through GitHub"
$d.Content.Find.Execute(
    "Using malicious scripts to uncover",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Using malicious scripts trawl through GitHub to uncover",
    2) | Out-Null

# Response strategy (col 3): append new sentence about rotating credentials
$d.Content.Find.Execute(
    "not using root or admin)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "not using root or admin) and change the credentials once every 3 months.",
    2) | Out-Null

# Likelihood (col 5): High -> Medium
$d.Tables.Item(1).Cell(2, 5).Range.Text = "Medium"

# --- Table 1, Row 3 (Repetitive Strain Injury) -----------------------------

# Risk Statement (col 2): shorten the ending of the sentence
$d.Content.Find.Execute(
    "resulting in injury and potential impact on further computer usage.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "resulting in RSI.",
    2) | Out-Null

# --- Table 1, Row 4 (SQL injections) ---------------------------------------

# Response strategy (col 3): reword ending of the sentence
$d.Content.Find.Execute(
    "Limit the number of characters able to be added into the input prompts, as well as prohibiting the use of special characters such as “;” etc.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Limit the number of characters able to be added into the input prompts as well as sanitizing the inputs.",
    2) | Out-Null

# Likelihood (col 5): High -> Medium
$d.Tables.Item(1).Cell(4, 5).Range.Text = "Medium"

# --- Table 2, Row 2 (Jenkins) -----------------------------------------------

# Response strategy (col 3): "version 2.133" -> "a recent full release"
$d.Content.Find.Execute(
    "version 2.133",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a recent full release",
    2) | Out-Null
